# "readme moved and bug report" -- add the new "buzzwire volume glitch" bug
# report row (row 9) to the Reports sheet, and update the sheet's
# selection/scroll position to match where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reports")

# --- New bug report row (row 9) -------------------------------------------
$ws.Range("B9").Value = "buzzwire volume glitch"
$ws.Range("D9").Value = "when the mouse is too far from the line the volume is set to 0, and when move to another visualser it stays at 0"
$ws.Range("E9").Value = "Yes"
$ws.Range("F9").Value = "command added to the select vis key press, sound.setVolume(vol);"
$ws.Range("G9").Value = "V0.5"
$ws.Range("H9").Value = "bug found by Deniz"

# Row grew to fit the wrapped text (two lines at the sheet's row height).
$ws.Rows.Item(9).RowHeight = 29

# --- View state left by the author after editing ---------------------------
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("H10").Select()
